$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts B:G to C:H)
$ws.Range("B:B").Insert()

# New header
$ws.Range("B1").Value = "deviceProfileId"

# New values for the inserted column
$ws.Range("B2").Value = "ddf28dd5-093e-4657-9a6d-b5b49904f942"
$ws.Range("B3").Value = "b0dac1ac-3234-44d3-b290-59b7f5803300"
$ws.Range("B4").Value = "ddf28dd5-093e-4657-9a6d-b5b49904f942"

# Adjust column width for the new column B (deviceProfileId / GUID values)
$ws.Range("B:B").ColumnWidth = 40

# Update selection to match the diff (active cell B8)
$ws.Range("B8").Select()
